$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("H2").Value = 2
$ws1.Range("L2").Value = 1.16

$ws1.Range("H3").Value = 1
$ws1.Range("J3").Value = "Normal"
$ws1.Range("L3").Value = 0.82

$ws1.Range("L4").Value = 0.82

$ws1.Range("L5").Value = 1.14

$ws1.Range("L6").Value = 0.9399999999999999

$ws1.Range("L7").Value = 1.05

$ws1.Range("L8").Value = 1.08

$ws1.Range("L9").Value = 0.88

$ws1.Range("L10").Value = 1.16

$ws1.Range("L11").Value = 0.88

$ws1.Range("L12").Value = 1.13

$ws1.Range("L13").Value = 0.84

$ws1.Range("L14").Value = 1.11

$ws1.Range("L15").Value = 0.96

$ws1.Range("L16").Value = 0.84

$ws1.Range("L17").Value = 1.05

# --- Sheet: Summary ---
# Values in column B are stored as text (e.g. "11", "6"), not numbers.
# Prefixing with an apostrophe forces Excel to keep them as text instead
# of auto-converting the numeric-looking strings into numbers.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'10"
$ws2.Range("B10").Value = "'5"
$ws2.Range("B11").Value = "'2"
$ws2.Range("B14").Value = "'0"
